# Generate Report for Archive
# Refreshes the localization-status report: files that were "Ready for
# handoff" and have since been picked up now show "In Translation",
# while files still awaiting handoff keep their "Ready for handoff"
# status.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (per-language status columns E = zh-cn, F = de-de) ---
$overview = $wb.Worksheets.Item("Overview")

# 5271b3f8-...md (row 3) and 9a1932f1-...md (row 4) have moved into translation
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

# ca7d6371-...md (row 5) is still awaiting handoff
$overview.Range("E5").Value = "Ready for handoff"
$overview.Range("F5").Value = "Ready for handoff"

# --- zh-cn sheet (Status column C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"
$zhcn.Range("C5").Value = "Ready for handoff"

# --- de-de sheet (Status column C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
$dede.Range("C5").Value = "Ready for handoff"
